$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of this block (rows 999 and 1000),
# shifting the existing rows 999-1071 down to 1001-1073.
$ws.Range("A999:A1000").EntireRow.Insert()

# Row 999
$ws.Cells.Item(999,1).Value2 = 10
$ws.Cells.Item(999,2).Value2 = 'Vega Modelo de Temuco'
$ws.Cells.Item(999,3).Value2 = 'La Araucanía'
$ws.Cells.Item(999,4).Value2 = 45021
$ws.Cells.Item(999,5).Value2 = 9
$ws.Cells.Item(999,6).Value2 = 100114001
$ws.Cells.Item(999,7).Value2 = 'Papa'
$ws.Cells.Item(999,8).Value2 = 'Rodeo'
$ws.Cells.Item(999,9).Value2 = '1a (guarda)'
$ws.Cells.Item(999,10).Value2 = 310
$ws.Cells.Item(999,11).Value2 = 12000
$ws.Cells.Item(999,12).Value2 = 12000
$ws.Cells.Item(999,13).Value2 = 12000
$ws.Cells.Item(999,14).Value2 = '$/malla 25 kilos'
$ws.Cells.Item(999,15).Value2 = 'Provincia de Cautín'
$ws.Cells.Item(999,16).Value2 = 480
$ws.Cells.Item(999,17).Value2 = 25
$ws.Cells.Item(999,18).Value2 = 'Hortaliza'

# Row 1000
$ws.Cells.Item(1000,1).Value2 = 10
$ws.Cells.Item(1000,2).Value2 = 'Vega Modelo de Temuco'
$ws.Cells.Item(1000,3).Value2 = 'La Araucanía'
$ws.Cells.Item(1000,4).Value2 = 45021
$ws.Cells.Item(1000,5).Value2 = 9
$ws.Cells.Item(1000,6).Value2 = 100114001
$ws.Cells.Item(1000,7).Value2 = 'Papa'
$ws.Cells.Item(1000,8).Value2 = 'Rosara'
$ws.Cells.Item(1000,9).Value2 = '1a (guarda)'
$ws.Cells.Item(1000,10).Value2 = 480
$ws.Cells.Item(1000,11).Value2 = 10000
$ws.Cells.Item(1000,12).Value2 = 10000
$ws.Cells.Item(1000,13).Value2 = 10000
$ws.Cells.Item(1000,14).Value2 = '$/saco 25 kilos'
$ws.Cells.Item(1000,15).Value2 = 'Provincia de Cautín'
$ws.Cells.Item(1000,16).Value2 = 400
$ws.Cells.Item(1000,17).Value2 = 25
$ws.Cells.Item(1000,18).Value2 = 'Hortaliza'
